$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 previously only had A18 (domain) and B18 = 0 (cookie count).
# Update it to reflect 4 cookies found, and list their names in C18:F18.
$ws.Range("B18").Value = 4
$ws.Range("C18").Value = "_gid"
$ws.Range("D18").Value = "_gat"
$ws.Range("E18").Value = "_ga"
$ws.Range("F18").Value = "NBGPUBLICConsent"
